# ------------------------------------------------------------------
# 1. Insert a new column at I (shifts old I->J, J->K) and update headers
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(9).Insert()

$ws.Range("I1").Value = "Searched Job Title"

# ------------------------------------------------------------------
# 2. Stamp the searched job title on every data row (new column I)
# ------------------------------------------------------------------
for ($r = 2; $r -le 65; $r++) {
    $ws.Cells.Item($r, 9).Value = "Python Developer"
}

# ------------------------------------------------------------------
# 3. Refresh rows 42-61 with the latest scraped listings (columns A-H)
# ------------------------------------------------------------------
$rowData = @{
    42 = @("Flutter Development", "Maxgen Technologies Private Limited", "Ahmedabad", "0-1", "₹ 5,000 - 12,000 /month", "1 week ago", "Internshala", "https://internshala.com/internship/detail/flutter-development-internship-in-ahmedabad-at-maxgen-technologies-private-limited1765353396")
    43 = @("MERN Stack", "Maxgen Technologies Private Limited", "Ahmedabad", "0-1", "₹ 5,000 - 12,000 /month", "1 week ago", "Internshala", "https://internshala.com/internship/detail/part-time-mern-stack-internship-in-ahmedabad-at-maxgen-technologies-private-limited1765357770")
    44 = @("Web Development", "Maxgen Technologies Private Limited", "Ahmedabad", "0-1", "₹ 5,000 - 12,000 /month", "2 weeks ago", "Internshala", "https://internshala.com/internship/detail/web-development-internship-in-ahmedabad-at-maxgen-technologies-private-limited1764577171")
    45 = @("Full Stack Development", "Maxgen Technologies Private Limited", "Ahmedabad", "0-1", "₹ 5,000 - 12,000 /month", "2 weeks ago", "Internshala", "https://internshala.com/internship/detail/full-stack-development-internship-in-ahmedabad-at-maxgen-technologies-private-limited1764824075")
    46 = @("Web Programmer", "Maxgen Technologies Private Limited", "Ahmedabad", "0-1", "₹ 5,000 - 12,000 /month", "3 days ago", "Internshala", "https://internshala.com/internship/detail/part-time-web-programmer-internship-in-ahmedabad-at-maxgen-technologies-private-limited1765780187")
    47 = @("Quality Analyst", "TSTEPS PRIVATE LIMITED", "Chennai, Mumbai, Bangalore, Kerala, Puduchery(Hybrid)", "0-1", "₹ 30,000 - 40,000 /month", "3 weeks ago", "Internshala", "https://internshala.com/internship/detail/quality-analyst-internship-in-multiple-locations-at-tsteps-private-limited1762161642")
    48 = @("Mobile App Development", "TSTEPS PRIVATE LIMITED", "Chennai, Mumbai, Hyderabad, Bangalore, Kerala, Puduchery(Hybrid)", "0-1", "₹ 30,000 - 40,000 /month", "3 weeks ago", "Internshala", "https://internshala.com/internship/detail/mobile-app-development-internship-in-multiple-locations-at-tsteps-private-limited1762159495")
    49 = @("Full Stack Development", "WonderBotz", "Ahmedabad", "0-1", "₹ 10,000 /month", "3 days ago", "Internshala", "https://internshala.com/internship/detail/full-stack-development-internship-in-ahmedabad-at-wonderbotz1765800276")
    50 = @("Full Stack Development", "Meru Technosoft Private Limited", "Ahmedabad", "0-1", "₹ 3,000 - 6,000 /month", "3 weeks ago", "Internshala", "https://internshala.com/internship/detail/full-stack-development-internship-in-ahmedabad-at-meru-technosoft-private-limited1763471030")
    51 = @("Product Management", "Coding Junior", "Ahmedabad, Bhubaneswar, Delhi, Surat, Hyderabad, Bhopal, Mumbai, Jaipur", "0-1", "₹ 15,000 /month", "2 weeks ago", "Internshala", "https://internshala.com/internship/detail/product-management-internship-in-multiple-locations-at-coding-junior1764562748")
    52 = @("ReactJS Development", "Chaintech Network", "Ahmedabad", "0-1", "₹ 9,500 - 10,000 /month", "2 weeks ago", "Internshala", "https://internshala.com/internship/detail/reactjs-development-internship-in-ahmedabad-at-chaintech-network1764573695")
    53 = @("Programming", "TSTEPS PRIVATE LIMITED", "Chennai, Coimbatore, Mumbai, Bangalore, Kerala(Hybrid)", "0-1", "₹ 30,000 - 40,000 /month", "3 weeks ago", "Internshala", "https://internshala.com/internship/detail/programming-internship-in-multiple-locations-at-tsteps-private-limited1762156258")
    54 = @("ReactJS Development", "Maxgen Technologies Private Limited", "Ahmedabad", "0-1", "₹ 5,000 - 12,000 /month", "2 weeks ago", "Internshala", "https://internshala.com/internship/detail/reactjs-development-internship-in-ahmedabad-at-maxgen-technologies-private-limited1764736218")
    55 = @("ReactJS Development", "Maxgen Technologies Private Limited", "Ahmedabad", "0-1", "₹ 5,000 - 12,000 /month", "1 week ago", "Internshala", "https://internshala.com/internship/detail/part-time-reactjs-development-internship-in-ahmedabad-at-maxgen-technologies-private-limited1764920480")
    56 = @("Data Anlayst", "Maxgen Technologies Private Limited", "Ahmedabad", "0-1", "₹ 5,000 - 12,000 /month", "1 week ago", "Internshala", "https://internshala.com/internship/detail/part-time-data-anlayst-internship-in-ahmedabad-at-maxgen-technologies-private-limited1765003178")
    57 = @("Embedded Systems", "TSTEPS PRIVATE LIMITED", "Chennai, Mumbai, Hyderabad, Bangalore, Kerala(Hybrid)", "0-1", "₹ 30,000 - 40,000 /month", "3 weeks ago", "Internshala", "https://internshala.com/internship/detail/embedded-systems-internship-in-multiple-locations-at-tsteps-private-limited1762156894")
    58 = @("JavaScript Development", "TSTEPS PRIVATE LIMITED", "Chennai, Coimbatore, Mumbai, Bangalore, Kerala(Hybrid)", "0-1", "₹ 30,000 - 40,000 /month", "3 weeks ago", "Internshala", "https://internshala.com/internship/detail/javascript-development-internship-in-multiple-locations-at-tsteps-private-limited1762158549")
    59 = @("Game Development", "TSTEPS PRIVATE LIMITED", "Chennai, Hyderabad, Mumbai, Kerala(Hybrid)", "0-1", "₹ 30,000 - 40,000 /month", "3 weeks ago", "Internshala", "https://internshala.com/internship/detail/game-development-internship-in-multiple-locations-at-tsteps-private-limited1762153391")
    60 = @("Software Testing", "TSTEPS PRIVATE LIMITED", "Chennai, Coimbatore, Mumbai, Hyderabad, Kerala(Hybrid)", "0-1", "₹ 30,000 - 40,000 /month", "3 weeks ago", "Internshala", "https://internshala.com/internship/detail/software-testing-internship-in-multiple-locations-at-tsteps-private-limited1762159771")
    61 = @("PHP Development", "ECodeSoft Solutions", "Ahmedabad", "0-1", "₹ 3,000 - 8,000 /month", "1 week ago", "Internshala", "https://internshala.com/internship/detail/php-development-internship-in-ahmedabad-at-ecodesoft-solutions1764927058")
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item([int]$r, $c + 1).Value = $vals[$c]
    }
}

Write-Output "edit complete"
